$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 8 & 9: height shrinks from 43.2 to 28.8
$ws.Rows.Item(8).RowHeight = 28.8
$ws.Rows.Item(9).RowHeight = 28.8

# Row 23: fill in the new journal entry (game input/position functions)
$ws.Range("B23").Value = 44261
$ws.Range("B22").Copy()
$ws.Range("B23").PasteSpecial(-4122)
$ws.Range("C23").Value = 0.75
$ws.Range("D23").Value = 0.77083333333333337
$ws.Range("E23").Formula = '=IF(ISBLANK(Tableau1[[#This Row],[Heure fin]]),"",Tableau1[[#This Row],[Heure fin]]-Tableau1[[#This Row],[Heure début]])'
$ws.Range("F23").Value = "Ma-20"
$ws.Range("G23").Value = "Code"
$ws.Range("H23").Value = "Jeu"
$ws.Range("I23").Value = "Maison"
$ws.Range("J23").Value = "j'ai crée 2 focntion pour demande ou le joueur voulais tirer et une troisieme qui regarde si il y a un bateau"
$ws.Range("K23").Value = "Oui"

$ws.Rows.Item(23).RowHeight = 57.6

# Update selection as left by the author
$ws.Range("L23").Select()
